# The workbook is a rail-car "trace report". This edit re-runs the trace
# search (new timestamp in the summary line) which pulls fresh event rows
# for the three non-BNGX cars, and additionally fixes row 7 so it carries
# the same highlight fill as the other data rows ("white color code" ->
# the olive/green row highlight used by rows 4-6).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 5: CGAX 10167 / WINDSOR -> JOHNSTOWN -------------------------
$ws.Range("A5").Value = "CGAX"
$ws.Range("B5").Value = 10167
$ws.Range("C5").Value = "WINDSOR"
$ws.Range("D5").Value = "CO"
$ws.Range("E5").Value = 6
$ws.Range("F5").Value = 16
$ws.Range("G5").Value = 1524
$ws.Range("H5").Value = "Arrive In-Transit"
$ws.Range("I5").ClearContents()
$ws.Range("J5").Value = "JOHNSTOWN"
$ws.Range("K5").Value = "CO"
$ws.Range("L5").Value = 273000
$ws.Range("M5").Value = 64200
$ws.Range("N5").Value = 208800
$ws.Range("O5").Value = "CGAX10167"

# ---- Row 6: CAIX 541012 / LA JUNTA -> LOVELAND ------------------------
$ws.Range("A6").Value = "CAIX"
$ws.Range("B6").Value = 541012
$ws.Range("C6").Value = "LA JUNTA"
$ws.Range("D6").Value = "CO"
$ws.Range("E6").Value = 6
$ws.Range("F6").Value = 22
$ws.Range("G6").Value = 719
$ws.Range("H6").Value = "Departure"
$ws.Range("I6").Value = "HKCKDE"
$ws.Range("J6").Value = "LOVELAND"
$ws.Range("K6").Value = "CO"
$ws.Range("L6").Value = 273100
$ws.Range("M6").Value = 62900
$ws.Range("N6").Value = 210200
$ws.Range("O6").Value = "CAIX541012"

# ---- Row 7: CGEX 1941 / GREELEY -> JOHNSTOWN --------------------------
$ws.Range("A7").Value = "CGEX"
$ws.Range("B7").Value = 1941
$ws.Range("C7").Value = "GREELEY"
$ws.Range("D7").Value = "CO"
$ws.Range("E7").Value = 6
$ws.Range("F7").Value = 21
$ws.Range("G7").Value = 1419
$ws.Range("H7").Value = "Junction Delivery"
$ws.Range("I7").Value = "GWR"
$ws.Range("J7").Value = "JOHNSTOWN"
$ws.Range("K7").Value = "CO"
$ws.Range("L7").Value = 198750
$ws.Range("M7").Value = 0
$ws.Range("N7").Value = 198750
$ws.Range("O7").Value = "CGEX1941"

# Row 7 had been left without the row highlight fill used elsewhere in the
# table (white/no fill) - apply the same olive fill (RGB 9F,A4,59, i.e.
# 159 + 164*256 + 89*65536 = 5874847) used by rows 4-6 to columns A:N (the
# Car_no column, O, is intentionally left unstyled, matching the rest of
# the table).
$ws.Range("A7:N7").Interior.Color = 5874847

# ---- Summary cells at the top of the report ---------------------------
$ws.Range("A1").Value = "Description unknown, completed 06/22/2023 11:08:20 EDT, by WPJTOWN1.The search returned: 4 events."
$ws.Range("A2").Value = "4 CO"
